# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Rows 45-51 also shuffled (new coin "Aave" inserted, later entries shifted down).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.572.17'
$ws.Range('E2').Value = '  +2.43%  '
# Row 3
$ws.Range('D3').Value = '1.682.94'
$ws.Range('E3').Value = '  +2.92%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.26%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.31'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.08%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5338'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.48%  '
# Row 7
$ws.Range('E7').Value = '  -0.29%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2680'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4.99%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06431'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.48%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.52'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +6.43%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07799'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.25%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.514'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.94%  '
# Row 13
$ws.Range('D13').Value = '1.667.78'
$ws.Range('E13').Value = '  +1.83%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5634'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.57%  '
# Row 15
$ws.Range('D15').Value = '0.0₅8439'
$ws.Range('E15').Value = '  +7.17%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.13'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.55%  '
# Row 17
$ws.Range('D17').Value = '26.584.76'
$ws.Range('E17').Value = '  +2.39%  '
# Row 18
$ws.Range('E18').Value = '  -0.18%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.800'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.82%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.68'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +6.00%  '
# Row 21
$ws.Range('E21').Value = '  +4.71%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.372'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.09%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.28%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.23'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.34%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1279'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +7.92%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.479'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.53%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.20'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.94%  '
# Row 28
$ws.Range('E28').Value = '  +3.26%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06146'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.66%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.278'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.97%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.608'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +8.43%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.463'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.85%  '
# Row 33
$ws.Range('E33').Value = '  +6.73%  '
# Row 34
$ws.Range('E34').Value = '  +5.25%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.418'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.53%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.788'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.01%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5716'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.76%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01647'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.56%  '
# Row 39
$ws.Range('E39').Value = '  +5.01%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8735'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.05%  '
# Row 41
$ws.Range('D41').Value = '1.062.71'
$ws.Range('E41').Value = '  +2.53%  '
# Row 42
$ws.Range('E42').Value = '  -0.09%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.07'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.61%  '
# Row 44
$ws.Range('D44').Value = '1.832.78'
$ws.Range('E44').Value = '  +2.42%  '
# Row 45
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.29'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.70%  '
# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.150'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.77%  '
# Row 47
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9992'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.16%  '
# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05205'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.38%  '
# Row 49
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.094'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.79%  '
# Row 50
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4254'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.51%  '
# Row 51
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09914'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.98%  '
